$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "DuplicateDeal_TC001"
$ws.Range("A4").VerticalAlignment = -4108
$ws.Range("B4").Value = "John Tucker"
$ws.Range("C4").Value = "ONE"
$ws.Range("D4").Value = "Deal shared successfully"

$ws.Range("A5").Value = "WithDrawDeal_TC001"
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("B5").Value = "John Tucker"
$ws.Range("C5").Value = "ONE"
$ws.Range("D5").Value = "Deal shared successfully"

$ws.Range("A6").Value = "WithDrawDeal_TC002"
$ws.Range("A6").VerticalAlignment = -4108
$ws.Range("B6").Value = "NA"
$ws.Range("C6").Value = "ALL"
$ws.Range("D6").Value = "Deal shared successfully"

$ws.Range("A7").Value = "WithDrawDeal_TC003"
$ws.Range("A7").VerticalAlignment = -4108
$ws.Range("B7").Value = "Stan Koster Andersons"
$ws.Range("C7").Value = "ONE"
$ws.Range("D7").Value = "Deal shared successfully"

[void]$ws.Range("B7").Select()
